# Apply "adding new progress as of date 04 nov 2025" update to the
# Training Dashboard sheet: for rows 3-32, decrement the value in column H
# (PERIOD TO EXPIRE) by 1 and update column I (LAST UPDATE) from
# "03-Nov-2025" to "04-Nov-2025", while preserving each cell's existing
# number format / style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 32; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE

    if ($iCell.Value2 -eq "03-Nov-2025") {
        # Decrement the numeric value in column H.
        $hCell.Value2 = $hCell.Value2 - 1

        # Writing a date-looking string via Value2 would normally get
        # auto-converted into a date serial number by Excel. Force the
        # cell to stay text, then restore the original (General) number
        # format by copying it back from the sibling H cell, which has
        # an identical style and was never touched by NumberFormat.
        $iCell.NumberFormat = "@"
        $iCell.Value2 = "04-Nov-2025"

        $hCell.Copy() | Out-Null
        $iCell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    }
}

$excel.CutCopyMode = 0
